$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily log entry (2026/01/25) was pushed in chronologically before the
# existing 2026/12/29 block, so insert a fresh row at 722 and shift the rest
# of the table (722-763) down by one to 723-764.
$ws.Rows.Item(722).Insert()

# Populate the newly inserted row with the pushed entry's data. Use a
# leading apostrophe to keep the date-looking text as a literal string
# (matching every other date cell in column A, which are plain text, not
# Excel date serials), then reset the style back to Normal so no stray
# quote-prefix style gets left behind on the cell.
$ws.Range("A722").Value = "'2026/01/25"
$ws.Range("A722").Style = "Normal"
$ws.Range("B722").Value = "日"
$ws.Range("C722").Value = 16
$ws.Range("D722").Value = 171
